$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking Price (D) and Volume(1h) (E) columns: force text with an
# apostrophe prefix so Excel stores them as literal strings (matches the source
# workbook, where these columns are inlineStr, not numbers/percentages). ---
$ws.Range("D2").Value = "'318.48"
$ws.Range("E2").Value = "'3.91%"
$ws.Range("D3").Value = "'39.60"
$ws.Range("E3").Value = "'2.09%"
$ws.Range("D4").Value = "'5.140"
$ws.Range("E4").Value = "'1.10%"
$ws.Range("D5").Value = "'0.08211"
$ws.Range("E5").Value = "'1.99%"
$ws.Range("D6").Value = "'2.007"
$ws.Range("E6").Value = "'4.41%"
$ws.Range("D7").Value = "'8.278"
$ws.Range("D8").Value = "'4.272"
$ws.Range("E8").Value = "'1.99%"
$ws.Range("D9").Value = "'0.9326"
$ws.Range("E9").Value = "'0.21%"
$ws.Range("D10").Value = "'0.1404"
$ws.Range("E10").Value = "'-3.27%"
$ws.Range("D11").Value = "'0.1987"
$ws.Range("E11").Value = "'2.54%"
$ws.Range("D12").Value = "'0.09035"
$ws.Range("E12").Value = "'0.56%"
$ws.Range("D13").Value = "'0.03584"
$ws.Range("E13").Value = "'2.41%"
$ws.Range("D14").Value = "'0.09813"
$ws.Range("E14").Value = "'0.14%"
$ws.Range("D15").Value = "'0.001401"
$ws.Range("E15").Value = "'0.83%"
$ws.Range("D16").Value = "'0.005990"
$ws.Range("E16").Value = "'-0.23%"
$ws.Range("D17").Value = "'3.667"
$ws.Range("E17").Value = "'-1.91%"
$ws.Range("D18").Value = "'3.175"
$ws.Range("E18").Value = "'-8.54%"
$ws.Range("D19").Value = "'0.3462"
$ws.Range("E19").Value = "'-0.03%"
$ws.Range("D20").Value = "'0.1276"
$ws.Range("E20").Value = "'-2.58%"
$ws.Range("D21").Value = "'4.899"
$ws.Range("E21").Value = "'1.96%"
$ws.Range("D22").Value = "'0.2451"
$ws.Range("E22").Value = "'2.00%"
$ws.Range("D23").Value = "'0.04335"
$ws.Range("E23").Value = "'-0.79%"
$ws.Range("D24").Value = "'0.001225"
$ws.Range("E24").Value = "'-0.62%"
$ws.Range("D25").Value = "'0.004776"
$ws.Range("E25").Value = "'11.60%"
$ws.Range("E26").Value = "'-0.03%"
$ws.Range("D27").Value = "'0.0004000"
$ws.Range("E27").Value = "'-10.05%"
$ws.Range("D39").Value = "'0.02215"
$ws.Range("E39").Value = "'6.76%"
$ws.Range("D40").Value = "'0.05267"
$ws.Range("E40").Value = "'4.18%"
$ws.Range("D41").Value = "'0.007504"
$ws.Range("E41").Value = "'0.90%"
$ws.Range("D42").Value = "'0.01011"
$ws.Range("E42").Value = "'0.09%"
$ws.Range("D43").Value = "'0.1380"
$ws.Range("E43").Value = "'2.26%"
$ws.Range("D44").Value = "'0.002149"
$ws.Range("E44").Value = "'0.43%"
$ws.Range("D45").Value = "'0.009869"
$ws.Range("E45").Value = "'10.40%"
$ws.Range("D46").Value = "'0.00006563"
$ws.Range("E46").Value = "'5.76%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.002767"
$ws.Range("E48").Value = "'-0.98%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.001200"
$ws.Range("E49").Value = "'-24.98%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'0.03%"
